# Update NATMI LR-pair output (Col8a1-Itga1) with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.274286666666667
$ws.Range("H2").Value = 6.82286
$ws.Range("I2").Value = 0.03663635825988129
$ws.Range("J2").Value = 0.03663635825988129
$ws.Range("M2").Value = 68.46613766666667
$ws.Range("N2").Value = 205.398413
$ws.Range("O2").Value = 0.4719163120948675
$ws.Range("P2").Value = 0.4719163120948675
$ws.Range("Q2").Value = 155.7116240134644
$ws.Range("R2").Value = 1401.40461612118
$ws.Range("S2").Value = 0.01728929507858951
$ws.Range("T2").Value = 0.01728929507858951

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.274286666666667
$ws.Range("H3").Value = 6.82286
$ws.Range("I3").Value = 0.03663635825988129
$ws.Range("J3").Value = 0.03663635825988129
$ws.Range("M3").Value = 9.278736333333333
$ws.Range("O3").Value = 0.06395551407683932
$ws.Range("P3").Value = 0.06395551407683933
$ws.Range("Q3").Value = 21.10250632641555
$ws.Range("R3").Value = 189.92255693774
$ws.Range("S3").Value = 0.002343097126413966
$ws.Range("T3").Value = 0.002343097126413967

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.274286666666667
$ws.Range("H4").Value = 6.82286
$ws.Range("I4").Value = 0.03663635825988129
$ws.Range("J4").Value = 0.03663635825988129
$ws.Range("M4").Value = 67.336226
$ws.Range("N4").Value = 202.008678
$ws.Range("O4").Value = 0.4641281738282933
$ws.Range("P4").Value = 0.4641281738282933
$ws.Range("Q4").Value = 153.1418809754533
$ws.Range("R4").Value = 1378.27692877908
$ws.Range("S4").Value = 0.01700396605487781
$ws.Range("T4").Value = 0.01700396605487781

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.9239867975814116
$ws.Range("J5").Value = 0.9239867975814117
$ws.Range("M5").Value = 68.46613766666667
$ws.Range("N5").Value = 205.398413
$ws.Range("O5").Value = 0.4719163120948675
$ws.Range("P5").Value = 0.4719163120948675
$ws.Range("Q5").Value = 3927.122990713653
$ws.Range("R5").Value = 35344.10691642288
$ws.Range("S5").Value = 0.4360444419389665
$ws.Range("T5").Value = 0.4360444419389666

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.9239867975814116
$ws.Range("J6").Value = 0.9239867975814117
$ws.Range("M6").Value = 9.278736333333333
$ws.Range("O6").Value = 0.06395551407683932
$ws.Range("P6").Value = 0.06395551407683933
$ws.Range("Q6").Value = 532.2154866805631
$ws.Range("R6").Value = 4789.939380125068
$ws.Range("S6").Value = 0.05909405063953165
$ws.Range("T6").Value = 0.05909405063953167

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.9239867975814116
$ws.Range("J7").Value = 0.9239867975814117
$ws.Range("M7").Value = 67.336226
$ws.Range("N7").Value = 202.008678
$ws.Range("O7").Value = 0.4641281738282933
$ws.Range("P7").Value = 0.4641281738282933
$ws.Range("Q7").Value = 3862.312819804851
$ws.Range("R7").Value = 34760.81537824366
$ws.Range("S7").Value = 0.4288483050029134
$ws.Range("T7").Value = 0.4288483050029135

# Row 8: MuSCs -> ECs
$ws.Range("G8").Value = 2.444408666666666
$ws.Range("H8").Value = 7.333226
$ws.Range("I8").Value = 0.03937684415870708
$ws.Range("J8").Value = 0.03937684415870709
$ws.Range("M8").Value = 68.46613766666667
$ws.Range("N8").Value = 205.398413
$ws.Range("O8").Value = 0.4719163120948675
$ws.Range("P8").Value = 0.4719163120948675
$ws.Range("Q8").Value = 167.3592202855931
$ws.Range("R8").Value = 1506.232982570338
$ws.Range("S8").Value = 0.01858257507731137
$ws.Range("T8").Value = 0.01858257507731137

# Row 9: MuSCs -> FAPs
$ws.Range("G9").Value = 2.444408666666666
$ws.Range("H9").Value = 7.333226
$ws.Range("I9").Value = 0.03937684415870708
$ws.Range("J9").Value = 0.03937684415870709
$ws.Range("M9").Value = 9.278736333333333
$ws.Range("O9").Value = 0.06395551407683932
$ws.Range("P9").Value = 0.06395551407683933
$ws.Range("Q9").Value = 22.68102350891489
$ws.Range("R9").Value = 204.129211580234
$ws.Range("S9").Value = 0.002518366310893699
$ws.Range("T9").Value = 0.0025183663108937

# Row 10: MuSCs -> MuSCs
$ws.Range("G10").Value = 2.444408666666666
$ws.Range("H10").Value = 7.333226
$ws.Range("I10").Value = 0.03937684415870708
$ws.Range("J10").Value = 0.03937684415870709
$ws.Range("M10").Value = 67.336226
$ws.Range("N10").Value = 202.008678
$ws.Range("O10").Value = 0.4641281738282933
$ws.Range("P10").Value = 0.4641281738282933
$ws.Range("Q10").Value = 164.5972544150253
$ws.Range("R10").Value = 1481.375289735228
$ws.Range("S10").Value = 0.01827590277050201
$ws.Range("T10").Value = 0.01827590277050202
